## [Closed] escape first n rows, sheet name or number #1
## Adds a new "Sheet2" (an options/reference sheet) after the existing
## "Sheet1", fills it in, formats the header rows, and moves the active
## selection on both sheets to match the post-edit workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- New worksheet, inserted directly after Sheet1 ------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# ---- Column widths (matches Sheet1's layout) -------------------------------
$ws2.Columns.Item(2).ColumnWidth = 18.625
$ws2.Columns.Item(3).ColumnWidth = 14.125

# ---- Row 1: title banner ----------------------------------------------------
$ws2.Range("A1").Value = "title"
[void]$ws2.Range("A1:C1").Merge()

# ---- Row 2: subtitle banner --------------------------------------------------
$ws2.Range("A2").Value = "subtitle1"
$ws2.Range("C2").Value = "subtitle2"
[void]$ws2.Range("A2:B2").Merge()

# Rows 1-2 (title/subtitle banner) are centered both ways; vertical center is
# already the sheet default, so only the horizontal alignment needs setting.
$ws2.Range("A1:C2").HorizontalAlignment = -4108

# ---- Row 3: column headers (same as Sheet1) ----------------------------------
$ws2.Range("A3").Value = "아이디"
$ws2.Range("B3").Value = "이름"
$ws2.Range("C3").Value = "이메일"

# ---- Rows 4-11: option rows --------------------------------------------------
$ws2.Range("A4").Value = 1
$ws2.Range("B4").Value = "1_option_name_xlsx"
$ws2.Range("C4").Value = "1_option_xlsx@email.com"

$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = "2_option_name_xlsx"
$ws2.Range("C5").Value = "2_option_xlsx@email.com"

$ws2.Range("A6").Value = 3
$ws2.Range("B6").Value = "3_option_name_xlsx"
$ws2.Range("C6").Value = "3_option_xlsx@email.com"

$ws2.Range("A7").Value = 4
$ws2.Range("C7").Value = "4_option_xlsx@email.com"

$ws2.Range("A8").Value = 5
$ws2.Range("B8").Value = "5_option_name_xlsx"
$ws2.Range("C8").Value = "5_option_xlsx@email.com"

$ws2.Range("A9").Value = 6
$ws2.Range("B9").Value = "6_option_name_xlsx"

$ws2.Range("A10").Value = 7
$ws2.Range("B10").Value = "7_option_name_xlsx"
$ws2.Range("C10").Value = "7_option_xlsx@email.com"

$ws2.Range("A11").Value = 8
$ws2.Range("B11").Value = "8_option_name_xlsx"
$ws2.Range("C11").Value = "8_option_xlsx@email.com"

# ---- Selections -------------------------------------------------------------
# Sheet1 loses the tab-selected flag and its cursor moves to C7.
[void]$ws1.Range("C7").Select()
# Sheet2 becomes the active (selected) tab, cursor resting on C9.
[void]$ws2.Range("C9").Select()
